$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "City" column header to "Country" (this also updates the
# bound table's column name automatically, since D3 is the table header
# cell for that column).
$ws.Range("D3").Value = "Country"

# Replace the city values with "India" for both data rows.
$ws.Range("D4").Value = "India"
$ws.Range("D5").Value = "India"

# Leave the selection on D5, matching the final state of the workbook.
$ws.Range("D5").Select()
